# Trading update: 2026-02-17 20:13:40
# Appends 4 new OPEN trades (trade #30-33) for the "MarketMaking" strategy
# to both the "All Trades" sheet and the dedicated "MarketMaking" sheet.

$wb = $excel.ActiveWorkbook

# Each entry: Num, Date, Time, Strategy, Side, Entry, Status, PnLPct, PnLUsd,
#             Capital, EntrySlip, ExitSlip, Confidence, Reason, Duration
$newTrades = @(
    @{ Num=30; Date="2026-02-17"; Time="20:11:43"; Strategy="MarketMaking"; Side="DOWN"; Entry=0.59; Status="OPEN"; PnLPct=0; PnLUsd=0; Capital=100; EntrySlip=0; ExitSlip=0; Confidence=0.6; Reason="Normal spread capture: 19600 bps"; Duration=0 },
    @{ Num=31; Date="2026-02-17"; Time="20:11:50"; Strategy="MarketMaking"; Side="UP";   Entry=0.41; Status="OPEN"; PnLPct=0; PnLUsd=0; Capital=100; EntrySlip=0; ExitSlip=0; Confidence=0.6; Reason="Normal spread capture: 19600 bps"; Duration=0 },
    @{ Num=32; Date="2026-02-17"; Time="20:11:57"; Strategy="MarketMaking"; Side="DOWN"; Entry=0.58; Status="OPEN"; PnLPct=0; PnLUsd=0; Capital=100; EntrySlip=0; ExitSlip=0; Confidence=0.6; Reason="Normal spread capture: 19600 bps"; Duration=0 },
    @{ Num=33; Date="2026-02-17"; Time="20:12:04"; Strategy="MarketMaking"; Side="DOWN"; Entry=0.53; Status="OPEN"; PnLPct=0; PnLUsd=0; Capital=100; EntrySlip=0; ExitSlip=0; Confidence=0.6; Reason="Normal spread capture: 19600 bps"; Duration=0 }
)

function Add-TradeRows {
    param($ws, $startRow, $trades)

    $r = $startRow
    foreach ($t in $trades) {
        $ws.Cells.Item($r, 1).Value  = $t.Num
        $ws.Cells.Item($r, 2).Value  = $t.Date      # fixed up to text below
        $ws.Cells.Item($r, 3).Value  = $t.Time
        $ws.Cells.Item($r, 4).Value  = $t.Strategy
        $ws.Cells.Item($r, 5).Value  = $t.Side
        $ws.Cells.Item($r, 6).Value  = $t.Entry
        # Column G (Exit Price) intentionally left blank (trade still OPEN)
        $ws.Cells.Item($r, 8).Value  = $t.Status
        $ws.Cells.Item($r, 9).Value  = $t.PnLPct
        $ws.Cells.Item($r, 10).Value = $t.PnLUsd
        $ws.Cells.Item($r, 11).Value = $t.Capital
        $ws.Cells.Item($r, 12).Value = $t.EntrySlip
        $ws.Cells.Item($r, 13).Value = $t.ExitSlip
        $ws.Cells.Item($r, 14).Value = $t.Confidence
        $ws.Cells.Item($r, 15).Value = $t.Reason
        # Column P (Exit Reason) intentionally left blank (trade still OPEN)
        $ws.Cells.Item($r, 17).Value = $t.Duration

        $r++
    }

    # The Date column ("2026-02-17") looks like a real date, and Excel's
    # auto-detection would otherwise silently convert it into a date serial
    # number. Force the column to Text first so the literal string is kept,
    # then restore the normal (unstyled) cell style so no stray number
    # format sticks to the newly written cells.
    $lastRow = $startRow + $trades.Count - 1
    $dateRange = $ws.Range("B" + $startRow + ":B" + $lastRow)
    $dateRange.NumberFormat = "@"

    $r = $startRow
    foreach ($t in $trades) {
        $ws.Cells.Item($r, 2).Value = $t.Date
        $r++
    }
    $dateRange.Style = "Normal"
}

# --- "All Trades" sheet: append after existing last row (30) ---
$wsAll = $wb.Worksheets.Item("All Trades")
Add-TradeRows $wsAll 31 $newTrades

# --- "MarketMaking" sheet: append after existing last row (2) ---
$wsMM = $wb.Worksheets.Item("MarketMaking")
Add-TradeRows $wsMM 3 $newTrades
